# final code with comment
# tut 6
#
# Update the attendance summary sheet:
#  - Column G ("Invalid") on row 3 becomes 1
#  - Column H ("Absent") on rows 3-18 becomes 1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G3: Invalid -> 1
$ws.Range("G3").Value = 1

# H3:H18: Absent -> 1 for every data row
$ws.Range("H3:H18").Value = 1
